$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-sample values for "albert" (row 11)
$ws.Range("G11").Value = 284910
$ws.Range("H11").Value = 140330

# Drop the border on G11:H11 (style index 1 -> 3, same font/fill/alignment, borderId 0)
$ws.Range("G11:H11").Borders.LineStyle = -4142

# Move the active selection
$ws.Range("G20").Select()
